# Insert a new data row at row 140 (shifting existing rows 140-239 down to 141-240)
# and populate it with the new "Choclo, Choclero, Primera" record dated 44673.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(140).Insert()

$ws.Range("A140").Value = 5
$ws.Range("B140").Value = "Macroferia Regional de Talca"
$ws.Range("C140").Value = "Maule"
$ws.Range("D140").Value = 44673
$ws.Range("E140").Value = 7
$ws.Range("F140").Value = 100112024
$ws.Range("G140").Value = "Choclo"
$ws.Range("H140").Value = "Choclero"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 40000
$ws.Range("K140").Value = 200
$ws.Range("L140").Value = 200
$ws.Range("M140").Value = 200
$ws.Range("N140").Value = "`$/unidad"
$ws.Range("O140").Value = "Región del Maule"
$ws.Range("P140").Value = 200
$ws.Range("Q140").Value = 1
$ws.Range("R140").Value = "Hortaliza"
